# Add team record columns (Wins / Losses / Ties) to the NYM_2022 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells go right after the existing last column (AC), i.e. AD, AE, AF.
# Copy formatting from an existing header cell (A1) so the new headers get the
# same bold font / thin border / centered alignment style already used by the
# other header cells (style index 1 in the original workbook).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Cells.Item(1, 30).Value = "Wins"
$ws.Cells.Item(1, 31).Value = "Losses"
$ws.Cells.Item(1, 32).Value = "Ties"

# Fill in the team record (same for every player row: 101-61-0) for all 62
# data rows (rows 2 through 63).
for ($r = 2; $r -le 63; $r++) {
    $ws.Cells.Item($r, 30).Value = 101
    $ws.Cells.Item($r, 31).Value = 61
    $ws.Cells.Item($r, 32).Value = 0
}
